$wb = $excel.ActiveWorkbook

# --- 1. commodity rename: elc_sol-CHE -> elc_spv-CHE (re_profiles column E, rows 4-387) ---
$wsRe = $wb.Worksheets.Item("re_profiles")
$wsRe.Range("E4:E387").Value = "elc_spv-CHE"

# --- 2. day/night timeslice membership lists re-shuffled (ev_charging_uc C13/C14) ---
$wsEv = $wb.Worksheets.Item("ev_charging_uc")
$wsEv.Range("C13").Value = "S2b0425h11,S2b0425h12,S2b0427h17,S3aH2,S4c1002h09,S4c1006h15,S4c1007h15,S2b0424h12,S2b0425h09,S2b0425h15,S2b0426h08,S2b0426h13,S2b0426h18,S2b0427h08,S2b0427h14,S2b0429h15,S2b0429h18,S4c1001h09,S4c1003h09,S4c1005h13,S4c1005h15,S4c1006h16,S4c1007h07,S4c1007h17,S6aH5,S1aH2,S1aH7,S2aH3,S2b0423h16,S2b0423h17,S2b0424h08,S2b0424h16,S2b0425h17,S2b0426h10,S2b0427h12,S2b0428h11,S3aH6,S4aH7,S4c1001h11,S4c1001h16,S4c1004h14,S4c1005h16,S4c1007h11,S5aH6,S5aH7,S6aH6,S2aH7,S2b0423h09,S2b0424h13,S2b0428h14,S3aH4,S4aH3,S4aH4,S4c1003h17,S2b0423h14,S2b0424h14,S2b0427h16,S2b0427h18,S2b0429h12,S3aH7,S4c1002h07,S4c1002h08,S4c1004h09,S4c1006h17,S4c1007h12,S1aH5,S2aH4,S2aH5,S2b0425h13,S2b0426h11,S2b0426h16,S2b0426h17,S2b0429h17,S4c1001h15,S4c1003h07,S4c1003h14,S4c1004h15,S4c1005h08,S4c1006h08,S4c1006h14,S4c1007h09,S4c1007h10,S5aH2,S5aH5,S6aH2,S1aH3,S2b0424h15,S2b0425h16,S2b0426h07,S2b0427h07,S2b0428h12,S2b0429h10,S4c1004h12,S4c1004h13,S4c1006h07,S2b0423h07,S2b0426h09,S2b0426h12,S2b0429h08,S4c1002h10,S4c1002h14,S4c1003h08,S4c1003h12,S4c1004h08,S4c1005h11,S6aH7,S2b0423h08,S2b0423h13,S2b0425h10,S2b0428h08,S2b0428h17,S3aH3,S4c1001h18,S4c1005h07,S4c1005h17,S1aH4,S2aH2,S2b0426h15,S2b0427h15,S2b0428h15,S4aH6,S4c1001h13,S4c1002h17,S4c1003h15,S4c1004h18,S2b0423h10,S2b0424h07,S2b0425h14,S2b0428h09,S4c1001h12,S4c1003h11,S4c1006h12,S4c1006h13,S4c1007h08,S4c1007h16,S2b0423h11,S2b0424h09,S2b0427h09,S2b0427h10,S2b0428h13,S2b0428h18,S4aH5,S4c1001h17,S4c1002h16,S4c1004h11,S4c1005h14,S4c1007h13,S4c1007h14,S2aH6,S2b0424h10,S2b0425h07,S2b0425h08,S2b0425h18,S2b0427h11,S2b0428h10,S2b0429h11,S2b0429h14,S3aH5,S4c1001h14,S4c1002h13,S4c1003h10,S4c1003h16,S4c1004h10,S4c1004h16,S4c1005h10,S4c1005h18,S4c1006h09,S4c1006h11,S5aH3,S2b0423h12,S2b0423h15,S2b0424h17,S2b0426h14,S2b0428h16,S2b0429h16,S4aH2,S4c1001h07,S4c1002h11,S4c1002h15,S4c1002h18,S4c1003h18,S4c1005h12,S4c1006h10,S4c1006h18,S4c1007h18,S2b0423h18,S2b0424h11,S2b0424h18,S2b0429h07,S4c1004h07,S5aH4,S6aH3,S6aH4,S1aH6,S2b0427h13,S2b0428h07,S2b0429h09,S2b0429h13,S4c1001h08,S4c1001h10,S4c1002h12,S4c1003h13,S4c1004h17,S4c1005h09"
$wsEv.Range("C14").Value = "S2b0423h04,S2b0424h21,S2b0428h06,S2b0429h02,S3aH8,S4c1001h01,S4c1001h03,S4c1002h22,S2b0423h22,S2b0426h05,S2b0429h23,S4aH8,S4c1001h19,S4c1004h04,S4c1005h05,S4c1006h23,S4c1007h21,S4c1007h24,S2b0423h01,S2b0424h04,S2b0425h20,S2b0425h22,S4c1002h06,S4c1002h24,S4c1003h01,S4c1004h03,S4c1004h19,S4c1004h24,S4c1006h01,S4c1006h04,S4c1007h04,S2b0423h20,S2b0424h23,S2b0425h04,S2b0426h04,S2b0427h04,S2b0427h23,S4c1003h05,S4c1003h23,S4c1004h05,S4c1005h20,S4c1005h23,S4c1006h05,S4c1006h06,S4c1006h20,S4c1007h03,S4c1007h20,S2aH8,S2b0424h05,S2b0424h20,S2b0424h22,S2b0425h06,S2b0425h23,S4c1001h23,S4c1002h04,S4c1003h21,S4c1004h20,S2b0423h23,S2b0424h03,S2b0426h20,S2b0427h20,S2b0428h05,S2b0429h21,S4c1002h19,S4c1002h20,S4c1003h04,S4c1005h02,S4c1005h24,S2b0424h19,S2b0426h21,S2b0427h03,S2b0428h04,S4c1002h02,S4c1003h06,S4c1004h06,S4c1007h19,S2b0423h06,S2b0424h02,S2b0424h24,S2b0426h03,S2b0428h22,S2b0429h04,S2b0429h05,S3aH1,S4c1003h03,S4c1005h21,S4c1006h03,S4c1006h22,S2b0424h06,S2b0425h01,S2b0425h05,S2b0425h21,S2b0428h02,S2b0428h20,S4c1003h02,S4c1004h23,S4c1006h21,S2b0423h21,S2b0427h02,S2b0427h06,S2b0429h03,S2b0429h20,S4aH1,S4c1003h20,S4c1004h21,S4c1006h19,S4c1007h06,S6aH1,S1aH1,S2b0425h03,S2b0426h23,S2b0427h05,S2b0427h21,S2b0427h22,S4c1001h02,S4c1001h05,S4c1001h06,S4c1004h02,S4c1005h22,S2b0426h02,S2b0428h03,S2b0428h23,S2b0428h24,S4c1002h05,S4c1003h24,S4c1007h22,S1aH8,S2b0423h03,S2b0423h19,S2b0423h24,S2b0425h02,S2b0425h24,S2b0426h24,S2b0427h24,S2b0428h01,S2b0428h19,S2b0428h21,S4c1001h22,S4c1002h03,S4c1003h19,S4c1007h01,S4c1007h05,S5aH1,S2aH1,S2b0423h02,S2b0423h05,S2b0426h01,S2b0426h06,S2b0426h19,S2b0427h01,S2b0429h06,S2b0429h19,S2b0429h22,S4c1004h22,S4c1005h19,S4c1007h02,S5aH8,S2b0424h01,S2b0426h22,S2b0427h19,S2b0429h01,S4c1001h04,S4c1001h20,S4c1002h01,S4c1005h04,S2b0425h19,S2b0429h24,S4c1001h21,S4c1001h24,S4c1002h21,S4c1002h23,S4c1003h22,S4c1004h01,S4c1005h01,S4c1005h03,S4c1005h06,S4c1006h02,S4c1006h24,S4c1007h23,S6aH8"

# --- 3. re_profiles: swap timeslice/ncap_afs pairs between rows 4<->9 and 5<->7 (M/N columns) ---
$m4 = $wsRe.Range("M4").Value()
$n4 = $wsRe.Range("N4").Value()
$m5 = $wsRe.Range("M5").Value()
$n5 = $wsRe.Range("N5").Value()
$m7 = $wsRe.Range("M7").Value()
$n7 = $wsRe.Range("N7").Value()
$m9 = $wsRe.Range("M9").Value()
$n9 = $wsRe.Range("N9").Value()

$wsRe.Range("M4").Value = $m9
$wsRe.Range("N4").Value = $n9
$wsRe.Range("M9").Value = $m4
$wsRe.Range("N9").Value = $n4

$wsRe.Range("M5").Value = $m7
$wsRe.Range("N5").Value = $n7
$wsRe.Range("M7").Value = $m5
$wsRe.Range("N7").Value = $n5
